$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.335.75"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.26%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.421.17"
$ws.Range("D3").ClearFormats()

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "669.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.01%  "

# Row 7
$ws.Range("E7").Value = "  -4.76%  "

# Row 8
$ws.Range("E8").Value = "  -4.79%  "

# Row 9
$ws.Range("E9").Value = "  -1.87%  "

# Row 10
$ws.Range("E10").Value = "  +0.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.418.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.35%  "

# Row 12
$ws.Range("E12").Value = "  +3.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +14.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "98.123.98"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.06%  "

# Row 16
$ws.Range("E16").Value = "  +0.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.053.38"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.11"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +20.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.591"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +34.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.421.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.44%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "512.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000208"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.55%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.601.38"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.153"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.71"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.62%  "

# Row 32
$ws.Range("E32").Value = "  +3.67%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.51"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +17.87%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.578"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "30.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.07"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +13.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "540.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.52%  "

# Row 41
$ws.Range("E41").Value = "  -2.90%  "

# Row 42
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.884"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.33%  "

# Row 44
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0440"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.52%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +14.48%  "

# Row 47
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.83"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +13.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +14.24%  "

# Row 50
$ws.Range("E50").Value = "  -2.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.98%  "

